$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Arkansas
$ws.Range("A7").Value = "Arkansas"
$ws.Range("B7").Value = "https://ardhslicensing.my.site.com/elicensing/s/search-provider/find-providers?language=en_US&tab=CC"
$ws.Range("C7").Value = "AR_childcare_providers_20250319.csv"
$ws.Range("D7").Value = ""
$ws.Range("D7").WrapText = $true

# Row 8: Alabama
$ws.Range("A8").Value = "Alabama"
$ws.Range("B8").Value = "https://apps.dhr.alabama.gov/daycare/daycare_search"

# Hyperlinks (Alabama's relationship is created before Arkansas's)
$ws.Hyperlinks.Add($ws.Range("B8"), "https://apps.dhr.alabama.gov/daycare/daycare_search")
$ws.Range("B8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ardhslicensing.my.site.com/elicensing/s/search-provider/find-providers?language=en_US&tab=CC")
$ws.Range("B7").Style = "Hyperlink"

# Selection matches target
$ws.Range("D7").Select()
